$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old layout (rows 1-2 = header over two rows, rows 3-14 = data):
#   Row1: E1 "mation"(0), G1 "pompes)"(1), I1 "Hiver"(2), J1 "Eté"(3), K1 "Année"(4)
#   Row2: F2 "(m3/s)"(5), G2 "(MW)"(6), H2 "(MW)"(6), I2 "(GWh)"(7), J2 "(GWh)"(7), K2 "(GWh)"(7)
#   Rows3-14: data
#
# New layout (single header row, data shifted up by one row):
#   Row1: A1 idx, B1 idx2, C1 Name, D1 Date Start, E1 Date End,
#         F1 (m3/s), G1 (MW1), H1 (MW2), I1 (GWh) Winter, J1 (GWh) Summer, K1 (GWh) Year
#   Rows2-13: same data as old rows3-14 (shifted up by one)

# Remove the old second header row - this shifts data rows 3..14 up to 2..13
$ws.Rows.Item(2).Delete()

# Rewrite row 1 completely with the new single-row header.
# Columns A-E: plain (no special style -> "Normal" clears any leftover style/format)
$ws.Cells.Item(1,1).Value = "idx"
$ws.Cells.Item(1,1).Style = "Normal"

$ws.Cells.Item(1,2).Value = "idx2"
$ws.Cells.Item(1,2).Style = "Normal"

$ws.Cells.Item(1,3).Value = "Name"
$ws.Cells.Item(1,3).Style = "Normal"

$ws.Cells.Item(1,4).Value = "Date Start"
$ws.Cells.Item(1,4).Style = "Normal"

$ws.Cells.Item(1,5).Value = "Date End"
$ws.Cells.Item(1,5).Style = "Normal"

# Columns F-K: keep the "header" font style (Arial 9, same as style used
# throughout the table for the text column) instead of the default.
$ws.Cells.Item(1,6).Value = "(m3/s)"
$ws.Cells.Item(1,6).Font.Name = "Arial"
$ws.Cells.Item(1,6).Font.Size = 9

$ws.Cells.Item(1,7).Value = "(MW1)"
$ws.Cells.Item(1,7).Font.Name = "Arial"
$ws.Cells.Item(1,7).Font.Size = 9

$ws.Cells.Item(1,8).Value = "(MW2)"
$ws.Cells.Item(1,8).Font.Name = "Arial"
$ws.Cells.Item(1,8).Font.Size = 9

$ws.Cells.Item(1,9).Value = "(GWh) Winter"
$ws.Cells.Item(1,9).Font.Name = "Arial"
$ws.Cells.Item(1,9).Font.Size = 9

$ws.Cells.Item(1,10).Value = "(GWh) Summer"
$ws.Cells.Item(1,10).Font.Name = "Arial"
$ws.Cells.Item(1,10).Font.Size = 9

$ws.Cells.Item(1,11).Value = "(GWh) Year"
$ws.Cells.Item(1,11).Font.Name = "Arial"
$ws.Cells.Item(1,11).Font.Size = 9

# Fix the selection to match the new intended view (row 2, the first data row).
$ws.Range("A2:K2").Select()
